$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 17; $row++) {
    $ws.Range("F$row").NumberFormat = "@"
    $ws.Range("F$row").Value = "1299"
    $ws.Range("K$row").Value = "Sediment"
    $ws.Range("L$row").Value = "Urban"
    $ws.Range("O$row").Value = "Flowing saltwater"
    $ws.Range("P$row").Value = "Harbour, marina scraped-off biofilm"
}
